$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# New row 8 values
$ws.Range("A8").Value = "CW3M"
$ws.Range("B8").Value = "Baseline 2010-18 C377+"
$ws.Range("C8").Value = "2010-18"

$ws.Range("D8").Value = 686.88716633333343
$ws.Range("E8").Value = 2094.2995878888887
$ws.Range("F8").Value = 5.8066811111111116
$ws.Range("G8").Value = 190.76777833333335
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 5.9919669999999989
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 535.46594233333337
$ws.Range("L8").Value = 89.974237444444441
$ws.Range("M8").Value = 1665.6726211111111
$ws.Range("N8").Value = 692.90829122222226
$ws.Range("O8").Value = 15597.417643111112
$ws.Range("P8").Value = 2216.8192002222222
$ws.Range("Q8").Value = 0.26791133333333339
$ws.Range("R8").Value = [double]"-1.1666666666666722E-5"

# Number formats to match the style of the row above (row 7), which uses
# a highlighted yellow fill on D, K, M/L, N, and a plain 0.00 / 0 / 0.000000
# format elsewhere.
$ws.Range("D8").NumberFormat = "0.00"
$ws.Range("D8").Interior.Color = 65535

$ws.Range("E8:J8").NumberFormat = "0.00"

$ws.Range("K8").NumberFormat = "0.00"
$ws.Range("K8").Interior.Color = 65535

$ws.Range("L8").NumberFormat = "0.00"

$ws.Range("M8").NumberFormat = "0.00"
$ws.Range("M8").Interior.Color = 65535

$ws.Range("N8").NumberFormat = "0.00"
$ws.Range("N8").Interior.Color = 65535

$ws.Range("O8:P8").NumberFormat = "0"
$ws.Range("Q8").NumberFormat = "0.00"
$ws.Range("R8").NumberFormat = "0.000000"

# Update selection to match the saved workbook state (B9 selected)
$ws.Range("B9").Select()
